# Applies the changes described by the commit:
#   "modifie excel with added channel and implement first slide 2 logic"
#
# - Rows 17-19: channel "CYW" -> new channel "PUE" (adds a new shared string)
#   and updated STATE/POS values for those rows (first "slide 2" logic).
# - Row 27 / Row 33: STATE value flips.
# - Rows 35-37: three new "CYW" rows appended at the end of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 17-19: channel CYW -> PUE, plus value tweaks ---

# Row 17: FIO CYW -> PUE, STATE 0 -> 1 (POS stays 1)
$ws.Cells.Item(17, 3).Value = "PUE"
$ws.Cells.Item(17, 4).Value = 1
$ws.Cells.Item(17, 5).Value = 1

# Row 18: FIO CYW -> PUE (POS/STATE unchanged: 2 / 1)
$ws.Cells.Item(18, 3).Value = "PUE"
$ws.Cells.Item(18, 4).Value = 2
$ws.Cells.Item(18, 5).Value = 1

# Row 19: FIO CYW -> PUE, POS 3 -> 1, STATE 1 -> 0
$ws.Cells.Item(19, 3).Value = "PUE"
$ws.Cells.Item(19, 4).Value = 1
$ws.Cells.Item(19, 5).Value = 0

# --- STATE flips further down the table ---

# Row 27: STATE 0 -> 1
$ws.Cells.Item(27, 5).Value = 1

# Row 33: STATE 1 -> 0
$ws.Cells.Item(33, 5).Value = 0

# --- Append three new CYW rows (35-37), continuing the Time sequence ---

$ws.Cells.Item(35, 1).Value = 17
$ws.Cells.Item(35, 2).Value = 34
$ws.Cells.Item(35, 3).Value = "CYW"
$ws.Cells.Item(35, 4).Value = 1
$ws.Cells.Item(35, 5).Value = 0

$ws.Cells.Item(36, 1).Value = 17
$ws.Cells.Item(36, 2).Value = 35
$ws.Cells.Item(36, 3).Value = "CYW"
$ws.Cells.Item(36, 4).Value = 2
$ws.Cells.Item(36, 5).Value = 1

$ws.Cells.Item(37, 1).Value = 17
$ws.Cells.Item(37, 2).Value = 36
$ws.Cells.Item(37, 3).Value = "CYW"
$ws.Cells.Item(37, 4).Value = 3
$ws.Cells.Item(37, 5).Value = 1

# --- Misc UI/state touch-ups reflected in the diff ---

# Selection moved to G17 in the saved view.
$ws.Range("G17").Select()

# Locale-style rename (Normal -> Standard) — best effort; harmless if the
# host does not persist it.
try {
    $wb.Styles.Item(1).Name = "Standard"
} catch {
}
